# Apply the "Manipulated data and working cat + budget class" edit:
# - Rename some existing categories
# - Add two new categories (Others, Rent) at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing categories
$ws.Range("B4").Value = "Investments"
$ws.Range("B5").Value = "Transportation"
$ws.Range("B6").Value = "Utilities"
$ws.Range("B7").Value = "Freelance"

# Append new categories
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Others"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Rent"
